# "Planned for today transactions"
#
# - mark "Today transactions" (Dashboard/Estimate row 97) and
#   "Base view" (Dashboard/Estimate row 117) as done (green, remaining
#   effort 0)
# - bump remaining effort for "Add tasks" (row 125) from 3 to 5
# - log five new backlog items on the Estimate sheet (rows 131-135)
# - log one new backlog idea on the Bugs sheet (row 9, Features(Global))
# - park the view on the newly added rows

$wb = $excel.ActiveWorkbook
$wsEstimate = $wb.Worksheets.Item("Estimate")
$wsBugs = $wb.Worksheets.Item("Bugs")

# --- Mark completed tasks as done: green fill + remaining hours -> 0 ---
# RGB(146,208,80) == the sheet's "done" green (matches the existing
# conditional fill already used for other completed rows).
$doneGreen = 146 + (208 * 256) + (80 * 65536)

$wsEstimate.Range("A97:D97").Interior.Color = $doneGreen
$wsEstimate.Range("D97").Value = 0

$wsEstimate.Range("A117:D117").Interior.Color = $doneGreen
$wsEstimate.Range("D117").Value = 0

# --- Bump remaining effort on an existing backlog row ---
$wsEstimate.Range("D125").Value = 5

# --- New backlog rows on the Estimate sheet ---
$wsEstimate.Range("C131").Value = "Cannot select empty category in category details"
$wsEstimate.Range("B131").Value = "Category details"
$wsEstimate.Range("A131").Value = "Bug"
$wsEstimate.Range("D131").Value = 1

$wsEstimate.Range("A132").Value = "Feature"
$wsEstimate.Range("B132").Value = "Reports page"
$wsEstimate.Range("C132").Value = "Show planned transactions"
$wsEstimate.Range("D132").Value = 2

$wsEstimate.Range("A133").Value = "Bug"
$wsEstimate.Range("B133").Value = "Record details"
$wsEstimate.Range("C133").Value = "Dublicate -> resultat amount without currency"
$wsEstimate.Range("D133").Value = 1

$wsEstimate.Range("A134").Value = "Bug"
$wsEstimate.Range("B134").Value = "Category selector"
$wsEstimate.Range("C134").Value = "Select new category"
$wsEstimate.Range("D134").Value = 1

# --- New backlog idea on the Bugs sheet (Features (Global) column) ---
$wsBugs.Range("D9").Value = "Filter templates (for reports)"

# --- Remaining new backlog row on the Estimate sheet ---
$wsEstimate.Range("A135").Value = "Bug"
$wsEstimate.Range("B135").Value = "Transactions list"
$wsEstimate.Range("C135").Value = "Update date record to last month -> exception"
$wsEstimate.Range("D135").Value = 1

# --- View state: leave the Bugs sheet parked near its new row, but make
#     sure Estimate stays the active tab/selection like before ---
$wsBugs.Range("A16").Select()

$wsEstimate.Activate()
$wsEstimate.Range("C137").Select()
